$p = $ppt.ActivePresentation

# -------------------------------------------------------------------
# Slide 8 ("Functional explanations" / Rules): fix stray space before
# colon in "...one of the following :" -> "...one of the following:"
# -------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shape8 = $s8.Shapes.Item(2)
$tr8 = $shape8.TextFrame.TextRange

for ($i = 1; $i -le $tr8.Paragraphs().Count; $i++) {
    $para = $tr8.Paragraphs($i)
    if ($para.Text -like "*can get killed if he encounters one of the following :*") {
        $ptext = $para.Text
        $idx = $ptext.IndexOf(" can get killed if he encounters one of the following :")
        $len = (" can get killed if he encounters one of the following :").Length
        $sub = $tr8.Characters($para.Start + $idx, $len)
        $sub.Text = " can get killed if he encounters one of the following:"
    }
}

# -------------------------------------------------------------------
# Slide 9 ("Functional explanations" / Functioning): bold key names
# and replace the placeholder "xxx" with "150".
# -------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shape9 = $s9.Shapes.Item(2)
$tr9 = $shape9.TextFrame.TextRange

# "The ZQSD keys" -> bold "ZQSD"
for ($i = 1; $i -le $tr9.Paragraphs().Count; $i++) {
    $para = $tr9.Paragraphs($i)
    if ($para.Text.TrimEnd() -eq "The ZQSD keys") {
        $idx = $para.Text.IndexOf("ZQSD")
        $sub = $tr9.Characters($para.Start + $idx, 4)
        $sub.Font.Bold = $true
    }
}

# "The WASD keys" -> bold "WASD"
for ($i = 1; $i -le $tr9.Paragraphs().Count; $i++) {
    $para = $tr9.Paragraphs($i)
    if ($para.Text.TrimEnd() -eq "The WASD keys") {
        $idx = $para.Text.IndexOf("WASD")
        $sub = $tr9.Characters($para.Start + $idx, 4)
        $sub.Font.Bold = $true
    }
}

# "...by the use of the SPACE key." -> bold "SPACE"
for ($i = 1; $i -le $tr9.Paragraphs().Count; $i++) {
    $para = $tr9.Paragraphs($i)
    if ($para.Text -like "*by the use of the SPACE key.*") {
        $idx = $para.Text.IndexOf("SPACE")
        $sub = $tr9.Characters($para.Start + $idx, 5)
        $sub.Font.Bold = $true
    }
}

# "...by the use of the R key. This action ... xxx point to the player."
for ($i = 1; $i -le $tr9.Paragraphs().Count; $i++) {
    $para = $tr9.Paragraphs($i)
    if ($para.Text -like "*by the use of the R key*") {
        # Bold the standalone "R" (the one immediately followed by " key.")
        $text = $para.Text
        $idxR = $text.IndexOf(" R key.")
        $rChar = $tr9.Characters($para.Start + $idxR + 1, 1)
        $rChar.Font.Bold = $true

        # Replace " xxx point to the " (incl. the yellow highlight run) with
        # plain " 150 point to the "
        $text2 = $para.Text
        $idxXxx = $text2.IndexOf(" xxx point to the ")
        $lenXxx = (" xxx point to the ").Length
        $subXxx = $tr9.Characters($para.Start + $idxXxx, $lenXxx)
        $subXxx.Text = " 150 point to the "
    }
}
